$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd shared string "lop-test" -> "log-test" in A2
$ws.Range("A2").Value = "log-test"

# Match the author's final active selection (A3)
$ws.Range("A3").Select() | Out-Null
